$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Modelo" header in F1, reusing the same formatting (border/bold/
# centered) already applied to the other header cells (e.g. E1).
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Cells.Item(1, 6).Value = "Modelo"

# Update the recomputed metric values for row 2.
$ws.Range("B2").Value = 0.501258630700655
$ws.Range("C2").Value = 0.9852767209371729
$ws.Range("D2").Value = 0.5170953347087761

# Add the model description / pipeline repr in F2 (spans two lines, as in
# the Python repr of the fitted sklearn Pipeline object).
$modelText = "Pipeline(steps=[('model',`n                 RandomForestRegressor(max_depth=7, n_estimators=150))])"
$ws.Range("F2").Value = $modelText

# Setting a value containing an embedded line break makes Excel mark the
# row as having a custom (auto-fitted) height; re-running AutoFit restores
# the row to its normal, non-"custom" state so the saved XML is unaffected.
$ws.Rows.Item(2).AutoFit()
